$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate and remove the rows for "Francisco Chahuán" and "Máximo Pacheco"
# (they dropped out of the candidate list). Delete the entire rows so the
# rest of the table shifts up, same as a manual row delete in Excel.
$namesToRemove = @("Francisco Chahuán", "Máximo Pacheco")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$rowsToDelete = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 2).Value()
    if ($namesToRemove -contains $name) {
        $rowsToDelete += $r
    }
}

# Delete from the bottom up so row indices of earlier rows stay valid.
$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Renumber the "id" column (A) so it stays a contiguous sequence after
# the deletions (same convention as the original sheet).
$newLastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $newLastRow; $r++) {
    $oldId = $ws.Cells.Item($r, 1).Value()
    $newId = $oldId
    foreach ($removedId in @(3, 19)) {
        if ($removedId -lt $oldId) {
            $newId = $newId - 1
        }
    }
    $ws.Cells.Item($r, 1).Value = $newId
}

$ws.Range("C14").Select()
